$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 132, pushing existing rows 132-159 down to 134-161.
$ws.Rows("132:133").Insert()

# New row 132 data
$ws.Cells.Item(132, 1).Value = 4
$ws.Cells.Item(132, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(132, 3).Value = "Los Lagos"
$ws.Cells.Item(132, 4).Value = 44551
$ws.Cells.Item(132, 5).Value = 10
$ws.Cells.Item(132, 6).Value = "Fruta"
$ws.Cells.Item(132, 7).Value = 100101
$ws.Cells.Item(132, 8).Value = "Berries"
$ws.Cells.Item(132, 9).Value = 100112025
$ws.Cells.Item(132, 10).Value = "Frutilla"
$ws.Cells.Item(132, 11).Value = "Sin especificar"
$ws.Cells.Item(132, 12).Value = "Primera"
$ws.Cells.Item(132, 13).Value = 1000
$ws.Cells.Item(132, 14).Value = 9500
$ws.Cells.Item(132, 15).Value = 10000
$ws.Cells.Item(132, 16).Value = 9750
$ws.Cells.Item(132, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(132, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(132, 19).Value = 1393
$ws.Cells.Item(132, 20).Value = 7

# New row 133 data
$ws.Cells.Item(133, 1).Value = 4
$ws.Cells.Item(133, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(133, 3).Value = "Los Lagos"
$ws.Cells.Item(133, 4).Value = 44551
$ws.Cells.Item(133, 5).Value = 10
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100101
$ws.Cells.Item(133, 8).Value = "Berries"
$ws.Cells.Item(133, 9).Value = 100112025
$ws.Cells.Item(133, 10).Value = "Frutilla"
$ws.Cells.Item(133, 11).Value = "Sin especificar"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 1000
$ws.Cells.Item(133, 14).Value = 9000
$ws.Cells.Item(133, 15).Value = 9500
$ws.Cells.Item(133, 16).Value = 9250
$ws.Cells.Item(133, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(133, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(133, 19).Value = 1321
$ws.Cells.Item(133, 20).Value = 7
